# Update cryptocurrency price/volume data (and fix row 43/44 ordering)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'41.529.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "
# Row 3
$ws.Range("D3").Value = "'2.478.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "
# Row 5
$ws.Range("D5").Value = "'313.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
# Row 6
$ws.Range("D6").Value = "'92.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.28%  "
# Row 7
$ws.Range("E7").Value = "  -0.66%  "
# Row 8
$ws.Range("E8").Value = "  -0.18%  "
# Row 9
$ws.Range("E9").Value = "  +1.52%  "
# Row 10
$ws.Range("E10").Value = "  -2.56%  "
# Row 11
$ws.Range("D11").Value = "'0.0787"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
# Row 12
$ws.Range("E12").Value = "  +2.22%  "
# Row 13
$ws.Range("D13").Value = "'2.859.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "
# Row 14
$ws.Range("D14").Value = "'6.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.19%  "
# Row 15
$ws.Range("D15").Value = "'16.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.24%  "
# Row 16
$ws.Range("D16").Value = "'2.513.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.45%  "
# Row 17
$ws.Range("E17").Value = "  -1.93%  "
# Row 18
$ws.Range("D18").Value = "'41.518.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.53%  "
# Row 19
$ws.Range("D19").Value = "'6.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.41%  "
# Row 20
$ws.Range("E20").Value = "  +2.31%  "
# Row 21
$ws.Range("D21").Value = "'72.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.41%  "
# Row 22
$ws.Range("D22").Value = "'11.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "
# Row 23
$ws.Range("D23").Value = "'235.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "
# Row 24
$ws.Range("E24").Value = "  -1.90%  "
# Row 25
$ws.Range("E25").Value = "  -0.09%  "
# Row 26
$ws.Range("E26").Value = "  -0.59%  "
# Row 27
$ws.Range("D27").Value = "'24.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "
# Row 28
$ws.Range("E28").Value = "  -0.03%  "
# Row 29
$ws.Range("D29").Value = "'9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "
# Row 30
$ws.Range("D30").Value = "'35.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.89%  "
# Row 31
$ws.Range("D31").Value = "'158.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.98%  "
# Row 32
$ws.Range("D32").Value = "'5.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.06%  "
# Row 33
$ws.Range("D33").Value = "'2.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "
# Row 34
$ws.Range("D34").Value = "'0.0754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.15%  "
# Row 35
$ws.Range("E35").Value = "  +2.75%  "
# Row 36
$ws.Range("E36").Value = "  -9.45%  "
# Row 37
$ws.Range("E37").Value = "  +3.83%  "
# Row 38
$ws.Range("E38").Value = "  -5.16%  "
# Row 39
$ws.Range("E39").Value = "  -3.62%  "
# Row 40
$ws.Range("E40").Value = "  +0.01%  "
# Row 41
$ws.Range("D41").Value = "'4.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.96%  "
# Row 42
$ws.Range("E42").Value = "  -0.21%  "
# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.963.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.74%  "
# Row 45
$ws.Range("D45").Value = "'0.0283"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.86%  "
# Row 46
$ws.Range("E46").Value = "  -3.37%  "
# Row 47
$ws.Range("D47").Value = "'8.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.34%  "
# Row 48
$ws.Range("D48").Value = "'2.717.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
# Row 49
$ws.Range("D49").Value = "'97.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.77%  "
# Row 50
$ws.Range("D50").Value = "'68.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.17%  "
# Row 51
$ws.Range("D51").Value = "'72.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.26%  "
